$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy style of A1 into B1:E1 and set header text
$ws.Range("A1").Copy() | Out-Null
$ws.Range("B1:E1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("B1").Value = "id_scooter.1"
$ws.Range("C1").Value = "id_scooter.2"
$ws.Range("D1").Value = "id_scooter.3"
$ws.Range("E1").Value = "id_scooter.4"

# Data rows 2-10, columns A-E
$data = @(
    @("300-200","300-200","300-200","300-200","300-200"),
    @("250-100","250-100","250-100","250-100","250-100"),
    @("123-423","123-423","123-423","123-423","123-423"),
    @("123-653","123-653","123-653","123-653","123-653"),
    @("566-124","566-124","566-124","566-124","566-124"),
    @("300-200","300-200","300-200","300-200","300-200"),
    @("250-100","250-100","250-100","250-100","250-100"),
    @("123-423","123-423","123-423","123-423","123-423"),
    @("123-653","123-653","123-653","123-653","123-653")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}
